$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 100.125
$ws.Range("I2").Value = 100.14286
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 100.14286
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = 12.85714
$ws.Range("N2").Value = -326

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4496.8057
$ws.Range("I32").Value = 3790.597
$ws.Range("J32").Value = 13960
$ws.Range("K32").Value = 3790.597
$ws.Range("L32").Value = 13960
$ws.Range("M32").Value = -3503.597
$ws.Range("N32").Value = -14534

$ws.Range("H74").Value = 4138.636
$ws.Range("I74").Value = 1045.6666
$ws.Range("J74").Value = 18057
$ws.Range("K74").Value = 1045.6666
$ws.Range("L74").Value = 18057
$ws.Range("M74").Value = -171.6666
$ws.Range("N74").Value = -19805

$ws.Range("H77").Value = 4138.636
$ws.Range("I77").Value = 1045.6666
$ws.Range("J77").Value = 18057
$ws.Range("K77").Value = 5228.333000000001
$ws.Range("L77").Value = 90285
$ws.Range("M77").Value = -860.3330000000005
$ws.Range("N77").Value = -99021

$ws.Range("H132").Value = 4040.5967
$ws.Range("I132").Value = 2624.9805
$ws.Range("J132").Value = 10603.909
$ws.Range("K132").Value = 7874.941500000001
$ws.Range("L132").Value = 31811.727
$ws.Range("M132").Value = -5344.941500000001
$ws.Range("N132").Value = -36871.727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 38219.965
$ws.Range("I20").Value = 771.1
$ws.Range("J20").Value = 145216.72
$ws.Range("K20").Value = 771.1
$ws.Range("L20").Value = 145216.72
$ws.Range("M20").Value = -524.1
$ws.Range("N20").Value = -145710.72

$ws.Range("H87").Value = 26500
$ws.Range("J87").Value = 26500
$ws.Range("L87").Value = 26500
$ws.Range("N87").Value = -28996

$ws.Range("H90").Value = 26500
$ws.Range("J90").Value = 26500
$ws.Range("L90").Value = 79500
$ws.Range("N90").Value = -91980

$ws.Range("H134").Value = 1434.6923
$ws.Range("I134").Value = 1142.5
$ws.Range("J134").Value = 1902.2
$ws.Range("K134").Value = 3427.5
$ws.Range("L134").Value = 5706.6
$ws.Range("M134").Value = -892.5
$ws.Range("N134").Value = -10776.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 953.4464
$ws.Range("I58").Value = 866.8889
$ws.Range("J58").Value = 1109.25
$ws.Range("K58").Value = 866.8889
$ws.Range("L58").Value = 1109.25
$ws.Range("M58").Value = -663.8889
$ws.Range("N58").Value = -1515.25

$ws.Range("H81").Value = 48000
$ws.Range("J81").Value = 48000
$ws.Range("L81").Value = 48000
$ws.Range("N81").Value = -49996

$ws.Range("H84").Value = 48000
$ws.Range("J84").Value = 48000
$ws.Range("L84").Value = 144000
$ws.Range("N84").Value = -153984

$ws.Range("H109").Value = 33000
$ws.Range("J109").Value = 33000
$ws.Range("L109").Value = 33000
$ws.Range("N109").Value = -35080

$ws.Range("H129:N129").ClearContents()

$ws.Range("H130:N130").ClearContents()

$ws.Range("H131:N131").ClearContents()

$ws.Range("H132:N132").ClearContents()

$ws.Range("H133:N133").ClearContents()

$ws.Range("H134:N134").ClearContents()

$ws.Range("H135:N135").ClearContents()

$ws.Range("H136").Value = 953.4464
$ws.Range("I136").Value = 866.8889
$ws.Range("J136").Value = 1109.25
$ws.Range("K136").Value = 2600.6667
$ws.Range("L136").Value = 3327.75
$ws.Range("M136").Value = -50.66670000000022
$ws.Range("N136").Value = -8427.75

$ws.Range("H137:N137").ClearContents()

$ws.Range("H138:N138").ClearContents()

$ws.Range("H139:N139").ClearContents()

$ws.Range("H140:N140").ClearContents()

$ws.Range("H141:N141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 461.7647
$ws.Range("I113").Value = 440
$ws.Range("J113").Value = 470.83334
$ws.Range("K113").Value = 1320
$ws.Range("L113").Value = 1412.50002
$ws.Range("M113").Value = 850
$ws.Range("N113").Value = -5752.500019999999

$ws.Range("H131").Value = 761.3913
$ws.Range("I131").Value = 506.26666
$ws.Range("J131").Value = 1239.75
$ws.Range("K131").Value = 1518.79998
$ws.Range("L131").Value = 3719.25
$ws.Range("M131").Value = 3521.20002
$ws.Range("N131").Value = -13799.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2238.125
$ws.Range("I68").Value = 2132.1428
$ws.Range("J68").Value = 2980
$ws.Range("K68").Value = 2132.1428
$ws.Range("L68").Value = 2980
$ws.Range("M68").Value = -1383.1428
$ws.Range("N68").Value = -4478

$ws.Range("H71").Value = 2238.125
$ws.Range("I71").Value = 2132.1428
$ws.Range("J71").Value = 2980
$ws.Range("K71").Value = 10660.714
$ws.Range("L71").Value = 14900
$ws.Range("M71").Value = -6916.714
$ws.Range("N71").Value = -22388

$ws.Range("H136").Value = 4494.5884
$ws.Range("I136").Value = 1552.64
$ws.Range("J136").Value = 12666.667
$ws.Range("K136").Value = 4657.92
$ws.Range("L136").Value = 38000.001
$ws.Range("M136").Value = -2107.92
$ws.Range("N136").Value = -43100.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 4143.9165
$ws.Range("I45").Value = 1784.5
$ws.Range("J45").Value = 4615.8
$ws.Range("K45").Value = 1784.5
$ws.Range("L45").Value = 4615.8
$ws.Range("M45").Value = -1293.5
$ws.Range("N45").Value = -5597.8
